$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.0006787744882400219
$ws.Range("J2").Value = 0.0006824388225951697
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1753453333333333
$ws.Range("N2").Value = 0.526036
$ws.Range("O2").Value = 0.414882210303281
$ws.Range("P2").Value = 0.5154067662594317
$ws.Range("Q2").Value = 0.0005997394884444444
$ws.Range("R2").Value = 0.005397655395999999
$ws.Range("S2").Value = 0.0002816114599784987
$ws.Range("T2").Value = 0.0003517335867236704

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.0006787744882400219
$ws.Range("J3").Value = 0.0006824388225951697
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2472935
$ws.Range("N3").Value = 0.494587
$ws.Range("O3").Value = 0.585117789696719
$ws.Range("P3").Value = 0.4845932337405682
$ws.Range("Q3").Value = 0.0008458262011666666
$ws.Range("R3").Value = 0.005074957206999999
$ws.Range("S3").Value = 0.0003971630282615232
$ws.Range("T3").Value = 0.0003307052358714992

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.954393
$ws.Range("H4").Value = 14.863179
$ws.Range("I4").Value = 0.9832128174003353
$ws.Range("J4").Value = 0.9885206487458582
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1753453333333333
$ws.Range("N4").Value = 0.526036
$ws.Range("O4").Value = 0.414882210303281
$ws.Range("P4").Value = 0.5154067662594317
$ws.Range("Q4").Value = 0.8687296920493333
$ws.Range("R4").Value = 7.818567228443999
$ws.Range("S4").Value = 0.4079175068815674
$ws.Range("T4").Value = 0.5094902309507784

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.954393
$ws.Range("H5").Value = 14.863179
$ws.Range("I5").Value = 0.9832128174003353
$ws.Range("J5").Value = 0.9885206487458582
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2472935
$ws.Range("N5").Value = 0.494587
$ws.Range("O5").Value = 0.585117789696719
$ws.Range("P5").Value = 0.4845932337405682
$ws.Range("Q5").Value = 1.2251891853455
$ws.Range("R5").Value = 7.351135112072999
$ws.Range("S5").Value = 0.5752953105187679
$ws.Range("T5").Value = 0.4790304177950798

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lgr6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.08117
$ws.Range("H6").Value = 0.16234
$ws.Range("I6").Value = 0.01610840811142459
$ws.Range("J6").Value = 0.01079691243154662
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1753453333333333
$ws.Range("N6").Value = 0.526036
$ws.Range("O6").Value = 0.414882210303281
$ws.Range("P6").Value = 0.5154067662594317
$ws.Range("Q6").Value = 0.01423278070666667
$ws.Range("R6").Value = 0.08539668423999999
$ws.Range("S6").Value = 0.006683091961735136
$ws.Range("T6").Value = 0.0055648017219297

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lgr6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.08117
$ws.Range("H7").Value = 0.16234
$ws.Range("I7").Value = 0.01610840811142459
$ws.Range("J7").Value = 0.01079691243154662
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2472935
$ws.Range("N7").Value = 0.494587
$ws.Range("O7").Value = 0.585117789696719
$ws.Range("P7").Value = 0.4845932337405682
$ws.Range("Q7").Value = 0.020072813395
$ws.Range("R7").Value = 0.08029125358
$ws.Range("S7").Value = 0.009425316149689458
$ws.Range("T7").Value = 0.005232110709616917
